# 9.c.1 workbook update:
#  - refresh the footnote text (new source agency: MDD/MCR KR instead of
#    SCITC/GKITS KR)
#  - append the 2023 data column (O) with its header + the three
#    technology rows' values, formatted like the existing 2022 column
#  - widen columns A:C slightly to fit the new footnote wording

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footnotes (B8/C8/A8 order keeps the shared-string table in the
#     same sequence the source edit produced) ---
$ws.Range("B8").Value = "*по данным МЦР КР"
$ws.Range("C8").Value = "*according to the MDD KR"
$ws.Range("A8").Value = "*КР СӨМ маалыматтары  боюнча"

# --- New 2023 column ---
$ws.Range("O4").Value = 2023
$ws.Range("O5").Value = 99
$ws.Range("O6").Value = 98.9
$ws.Range("O7").Value = 98.8

# Match formatting of the preceding (2022 / column N) cells.
$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial(-4122)
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("N6").Copy()
$ws.Range("O6").PasteSpecial(-4122)
$ws.Range("N7").Copy()
$ws.Range("O7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# PasteSpecial(xlPasteFormats) shouldn't disturb the values already
# written above, but re-assert them so formatting-only paste ordering
# can never clobber the data.
$ws.Range("O4").Value = 2023
$ws.Range("O5").Value = 99
$ws.Range("O6").Value = 98.9
$ws.Range("O7").Value = 98.8

# --- Column widths ---
$ws.Range("A1:C1").ColumnWidth = 37.140625
